# Rebuild rows 2-19 of Sheet1 to match the updated TPM-based NATMI export
# (Sending cluster now iterates ECs, FAPs, MuSCs; ligand/receptor stay Bmp7/Acvr2a;
#  Target cluster iterates ECs, FAPs, Inflammatory-Mac, MuSCs, Neutrophils, Resolving-Mac)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,("ECs", "Bmp7", "Acvr2a", "ECs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 11.651608, 34.954824, 0.1892813629236475, 0.1892813629236474, 1.516227632909333, 13.646048696184, 0.05652181670932983, 0.05652181670932983)
    ,("ECs", "Bmp7", "Acvr2a", "FAPs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 27.20435933333333, 81.613078, 0.4419371310876561, 0.4419371310876561, 3.540112348166444, 31.861011133498, 0.1319680349642224, 0.1319680349642224)
    ,("ECs", "Bmp7", "Acvr2a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 8.657179, 25.971537, 0.1406366091439035, 0.1406366091439035, 1.126561588996333, 10.139054300967, 0.04199587599049498, 0.04199587599049498)
    ,("ECs", "Bmp7", "Acvr2a", "MuSCs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 5.488499666666667, 16.465499, 0.08916114387925267, 0.08916114387925267, 0.7142202911232222, 6.427982620109001, 0.0266246489041299, 0.0266246489041299)
    ,("ECs", "Bmp7", "Acvr2a", "Neutrophils", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 4.091608333333333, 12.274825, 0.06646852536431769, 0.06646852536431769, 0.532442356286111, 4.791981206575, 0.01984834507503455, 0.01984834507503455)
    ,("ECs", "Bmp7", "Acvr2a", "Resolving-Mac", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.2986126887311651, 0.2986126887311651, 3, 1, 4.463825666666667, 13.391477, 0.07251522760122259, 0.07251522760122257, 0.5808791219452222, 5.227912097507, 0.02165396708795347, 0.02165396708795347)
    ,("FAPs", "Bmp7", "Acvr2a", "ECs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 11.651608, 34.954824, 0.1892813629236475, 0.1892813629236474, 3.311604490282667, 29.804440412544, 0.123449736669417, 0.1234497366694169)
    ,("FAPs", "Bmp7", "Acvr2a", "FAPs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 27.20435933333333, 81.613078, 0.4419371310876561, 0.4419371310876561, 7.731986737240889, 69.587880635168, 0.2882324050002536, 0.2882324050002535)
    ,("FAPs", "Bmp7", "Acvr2a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 8.657179, 25.971537, 0.1406366091439035, 0.1406366091439035, 2.460531872474667, 22.144786852272, 0.09172351729048953, 0.0917235172904895)
    ,("FAPs", "Bmp7", "Acvr2a", "MuSCs", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 5.488499666666667, 16.465499, 0.08916114387925267, 0.08916114387925267, 1.559934057260445, 14.039406515344, 0.05815110142395648, 0.05815110142395646)
    ,("FAPs", "Bmp7", "Acvr2a", "Neutrophils", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 4.091608333333333, 12.274825, 0.06646852536431769, 0.06646852536431769, 1.162911465022222, 10.4662031852, 0.04335092386427623, 0.04335092386427622)
    ,("FAPs", "Bmp7", "Acvr2a", "Resolving-Mac", 1, 0.3333333333333333, 0.2842186666666667, 0.852656, 0.6522022811047395, 0.6522022811047393, 3, 1, 4.463825666666667, 13.391477, 0.07251522760122259, 0.07251522760122257, 1.268702579212444, 11.418323212912, 0.04729459685634674, 0.04729459685634672)
    ,("MuSCs", "Bmp7", "Acvr2a", "ECs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 11.651608, 34.954824, 0.1892813629236475, 0.1892813629236474, 0.249740565872, 2.247665092848, 0.009309809544900694, 0.009309809544900694)
    ,("MuSCs", "Bmp7", "Acvr2a", "FAPs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 27.20435933333333, 81.613078, 0.4419371310876561, 0.4419371310876561, 0.5830982379506666, 5.247884141556, 0.02173669112318016, 0.02173669112318016)
    ,("MuSCs", "Bmp7", "Acvr2a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 8.657179, 25.971537, 0.1406366091439035, 0.1406366091439035, 0.185557974686, 1.670021772174, 0.006917215862918992, 0.006917215862918992)
    ,("MuSCs", "Bmp7", "Acvr2a", "MuSCs", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 5.488499666666667, 16.465499, 0.08916114387925267, 0.08916114387925267, 0.1176405018553333, 1.058764516698, 0.004385393551166295, 0.004385393551166295)
    ,("MuSCs", "Bmp7", "Acvr2a", "Neutrophils", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 4.091608333333333, 12.274825, 0.06646852536431769, 0.06646852536431769, 0.08769953301666666, 0.78929579715, 0.003269256425006908, 0.003269256425006908)
    ,("MuSCs", "Bmp7", "Acvr2a", "Resolving-Mac", 1, 0.3333333333333333, 0.021434, 0.064302, 0.04918503016409543, 0.04918503016409543, 3, 1, 4.463825666666667, 13.391477, 0.07251522760122259, 0.07251522760122257, 0.09567763933933332, 0.8610987540539999, 0.003566663656922378, 0.003566663656922378)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
